$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data (2021-09-10) is inserted above the existing
# "Acelga" rows, pushing the rest of the table down by two rows.
$ws.Rows.Item(194).Insert()
$ws.Rows.Item(194).Insert()

# Row 194 - "Primera" quality
$ws.Cells.Item(194, 1).Value = 8
$ws.Cells.Item(194, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(194, 3).Value = "Coquimbo"
$ws.Cells.Item(194, 4).Value = 44449
$ws.Cells.Item(194, 5).Value = 4
$ws.Cells.Item(194, 6).Value = 100112009
$ws.Cells.Item(194, 7).Value = "Acelga"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 3500
$ws.Cells.Item(194, 11).Value = 450
$ws.Cells.Item(194, 12).Value = 500
$ws.Cells.Item(194, 13).Value = 475
$ws.Cells.Item(194, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(194, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(194, 16).Value = 238
$ws.Cells.Item(194, 17).Value = 2
$ws.Cells.Item(194, 18).Value = "Hortaliza"

# Row 195 - "Segunda" quality
$ws.Cells.Item(195, 1).Value = 8
$ws.Cells.Item(195, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(195, 3).Value = "Coquimbo"
$ws.Cells.Item(195, 4).Value = 44449
$ws.Cells.Item(195, 5).Value = 4
$ws.Cells.Item(195, 6).Value = 100112009
$ws.Cells.Item(195, 7).Value = "Acelga"
$ws.Cells.Item(195, 8).Value = "Sin especificar"
$ws.Cells.Item(195, 9).Value = "Segunda"
$ws.Cells.Item(195, 10).Value = 1600
$ws.Cells.Item(195, 11).Value = 350
$ws.Cells.Item(195, 12).Value = 400
$ws.Cells.Item(195, 13).Value = 375
$ws.Cells.Item(195, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(195, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(195, 16).Value = 188
$ws.Cells.Item(195, 17).Value = 2
$ws.Cells.Item(195, 18).Value = "Hortaliza"
